$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T for each data row (row 1 is the header and is unchanged).
# This expands the sheet from 6 sender/target pairs (rows 2-7) into the
# full 3x3 sender x target cluster matrix (rows 2-10), per "Natmi following
# Dr Hou advice".
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "Rarres2", "Gpr1", "ECs", 3, 1, 1.248266, 3.744798, 0.008714053836811378, 0.00871405383681138, 3, 1, 1.446520666666667, 4.339562, 0.2269546876364677, 0.2269546876364677, 1.805642566497333, 16.250783098476, 0.001977695366580889, 0.00197769536658089),
    @("ECs", "Rarres2", "Gpr1", "FAPs", 3, 1, 1.248266, 3.744798, 0.008714053836811378, 0.00871405383681138, 3, 1, 4.543978333333333, 13.631935, 0.7129363631181282, 0.7129363631181282, 5.672093658236666, 51.04884292413, 0.006212565850431876, 0.006212565850431877),
    @("ECs", "Rarres2", "Gpr1", "sCs", 3, 1, 1.248266, 3.744798, 0.008714053836811378, 0.00871405383681138, 2, 0.6666666666666666, 0.383111, 1.149333, 0.0601089492454041, 0.0601089492454041, 0.478224435526, 4.304019919734, 0.000523792619798614, 0.0005237926197986142),
    @("FAPs", "Rarres2", "Gpr1", "ECs", 3, 1, 103.907654, 311.722962, 0.7253717484997341, 0.7253717484997342, 3, 1, 1.446520666666667, 4.339562, 0.2269546876364677, 0.2269546876364677, 150.3045689358493, 1352.741120422644, 0.1646265186010755, 0.1646265186010755),
    @("FAPs", "Rarres2", "Gpr1", "FAPs", 3, 1, 103.907654, 311.722962, 0.7253717484997341, 0.7253717484997342, 3, 1, 4.543978333333333, 13.631935, 0.7129363631181282, 0.7129363631181282, 472.1541284434966, 4249.38715599147, 0.517143896284038, 0.5171438962840381),
    @("FAPs", "Rarres2", "Gpr1", "sCs", 3, 1, 103.907654, 311.722962, 0.7253717484997341, 0.7253717484997342, 2, 0.6666666666666666, 0.383111, 1.149333, 0.0601089492454041, 0.0601089492454041, 39.808165231594, 358.273487084346, 0.04360133361462054, 0.04360133361462055),
    @("sCs", "Rarres2", "Gpr1", "ECs", 3, 1, 38.091531, 114.274593, 0.2659141976634544, 0.2659141976634544, 3, 1, 1.446520666666667, 4.339562, 0.2269546876364677, 0.2269546876364677, 55.100186816474, 495.901681348266, 0.06035047366881122, 0.06035047366881122),
    @("sCs", "Rarres2", "Gpr1", "FAPs", 3, 1, 38.091531, 114.274593, 0.2659141976634544, 0.2659141976634544, 3, 1, 4.543978333333333, 13.631935, 0.7129363631181282, 0.7129363631181282, 173.087091547495, 1557.783823927455, 0.1895799009836583, 0.1895799009836583),
    @("sCs", "Rarres2", "Gpr1", "sCs", 3, 1, 38.091531, 114.274593, 0.2659141976634544, 0.2659141976634544, 2, 0.6666666666666666, 0.383111, 1.149333, 0.0601089492454041, 0.0601089492454041, 14.593284532941, 131.339560796469, 0.01598382301098494, 0.01598382301098494)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $ws.Range("$($columns[$c])$r").Value = $rowData[$c]
    }
}
